$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal (non-date-coerced) text value into a cell.
#
# Assigning a date-shaped string straight to Range.Value (e.g. "2024-10-26")
# gets auto-converted into a date serial by Excel's usual type inference, the
# same way it would in real Excel. The source data stores these as plain
# text, so we round-trip the text through a throwaway formula cell instead:
# "=""literal""" always evaluates to a string, and pasting *values only*
# from it into the destination keeps that string typing without touching the
# destination's number format / style.
# ---------------------------------------------------------------------------
function Set-LiteralText {
    param($ws, [string]$cellRef, [string]$text)
    $scratch = $ws.Range("ZZ999")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Worksheets.Item(1))
#   - a handful of "想去人数" (F column) counters went up
#   - a brand-new event ("合肥·W·A第五人格同人only2.0", 2024-10-26) was appended
#     as the new last row (row 10)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 180
$ws1.Range("F3").Value = 497
$ws1.Range("F4").Value = 22
$ws1.Range("F7").Value = 31
$ws1.Range("F8").Value = 21
$ws1.Range("F9").Value = 1658

# Give the new row's A cell the same style (bold / bordered / centered) as
# the existing index column cells, by copying format from the row above.
$ws1.Range("A9").Copy()
$ws1.Range("A10").PasteSpecial(-4122)

$ws1.Range("A10").Value = 9
Set-LiteralText $ws1 "B10" "2024-10-26"
$ws1.Range("C10").Value = "合肥·W·A第五人格同人only2.0"
$ws1.Range("D10").Value = "莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)"
$ws1.Range("E10").Value = "2024.10.26 09:30-10.26 17:00"
$ws1.Range("F10").Value = 1
$ws1.Range("G10").Value = 68
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=91123"
$ws1.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png"

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (Worksheets.Item(4))
#   - same kind of F-column bumps
#   - the new "W·A第五人格" event also needs to land here, in date order; it
#     sorts before the already-present 2024-10-26 19:30 "钢琴小提琴" event, so
#     that event (and everything after it) shifts down by one row
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 180
$ws4.Range("F4").Value = 497
$ws4.Range("F5").Value = 22
$ws4.Range("F8").Value = 31
$ws4.Range("F9").Value = 21
$ws4.Range("F10").Value = 1658

# Give row 13 (brand new) the same column-A style as the rows above it.
$ws4.Range("A12").Copy()
$ws4.Range("A13").PasteSpecial(-4122)

# Push the old row 12 ("谭小棠独奏音乐会") down to row 13, renumbering its
# index in column A (11 -> 12).
$ws4.Range("A13").Value = 12
Set-LiteralText $ws4 "B13" $ws4.Range("B12").Value2
$ws4.Range("C13").Value = $ws4.Range("C12").Value2
$ws4.Range("D13").Value = $ws4.Range("D12").Value2
$ws4.Range("E13").Value = $ws4.Range("E12").Value2
$ws4.Range("F13").Value = $ws4.Range("F12").Value2
$ws4.Range("G13").Value = $ws4.Range("G12").Value2
$ws4.Range("H13").Value = $ws4.Range("H12").Value2
$ws4.Range("I13").Value = $ws4.Range("I12").Value2

# Push the old row 11 ("钢琴小提琴唯美经典音乐集") down to row 12, renumbering
# its index in column A (10 -> 11).
$ws4.Range("A12").Value = 11
Set-LiteralText $ws4 "B12" $ws4.Range("B11").Value2
$ws4.Range("C12").Value = $ws4.Range("C11").Value2
$ws4.Range("D12").Value = $ws4.Range("D11").Value2
$ws4.Range("E12").Value = $ws4.Range("E11").Value2
$ws4.Range("F12").Value = $ws4.Range("F11").Value2
$ws4.Range("G12").Value = $ws4.Range("G11").Value2
$ws4.Range("H12").Value = $ws4.Range("H11").Value2
$ws4.Range("I12").Value = $ws4.Range("I11").Value2

# Overwrite row 11 with the new "W·A第五人格" event. A11/B11 already hold the
# right values (10 / 2024-10-26) so only C..I need to change.
$ws4.Range("C11").Value = "合肥·W·A第五人格同人only2.0"
$ws4.Range("D11").Value = "莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)"
$ws4.Range("E11").Value = "2024.10.26 09:30-10.26 17:00"
$ws4.Range("F11").Value = 1
$ws4.Range("G11").Value = 68
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=91123"
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png"
